# CertificateNR12.pptx — "remove RG to user and models certificate"
#
# The certificate paragraph used to read:
#   "...portador do RG nº {{RG}} e CPF nº {{CPF}}..."
# The RG reference is dropped, leaving only the CPF reference:
#   "...portador do CPF nº {{CPF}}..."
# The shape is also shortened (height) to tighten up the now-shorter text block.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Rectangle 5")
$tr = $sh.TextFrame.TextRange

# --- 1. Rewrite "portador do RG n\u00ba {{RG}} e CPF n\u00ba " -> "portador do CPF n\u00ba " ---
# (the following "{{CPF}}" run is left completely untouched)
$para1 = $tr.Paragraphs(1, 1)

$oldChunk = "portador do RG n" + [char]0x00BA + " {{RG}} e CPF n" + [char]0x00BA + " "
$newChunk = "portador do CPF n" + [char]0x00BA + " "

$target = $para1.Characters(28, $oldChunk.Length)
if ($target.Text -ne $oldChunk) {
    throw "unexpected source text: [$($target.Text)]"
}
$target.Text = $newChunk

# --- 2. Split the single new run back into four runs, matching how the
#        author actually retyped the sentence word by word ("portador " /
#        "do " / "CPF " / "nº ") instead of leaving one big run behind. ---
$para1 = $tr.Paragraphs(1, 1)
$doRun = $para1.Characters(37, 3)          # "do "
$doRun.Text = "do "

$para1 = $tr.Paragraphs(1, 1)
$cpfRun = $para1.Characters(40, 4)         # "CPF "
$cpfRun.Text = "CPF "

$para1 = $tr.Paragraphs(1, 1)
$noRun = $para1.Characters(44, 3)          # "nº "
$noRun.Text = "nº "

# --- 3. Shrink the textbox now that a whole clause has been removed ---
$sh.Height = 2523768 / 12700
